$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 5 (shifts existing rows 5-22 down to 6-23)
$ws.Rows.Item(5).Insert()
# The inserted row inherits column A's style, which would materialize an
# empty-but-styled cell; clear it so the row is truly empty (absent from sheetData).
$ws.Range("A5").Clear()

# Add new column C header and "Cosmed" device values for the first three data rows
$ws.Range("C1").Value = "Device"
$ws.Range("C2").Value = "Cosmed"
$ws.Range("C3").Value = "Cosmed"
$ws.Range("C4").Value = "Cosmed"

# Update selection to match target (whole row 5 selected)
$ws.Range("A5:XFD5").Select()
